# Update countries & provincias Spain
# - Swap Iran / Estados Unidos positions (rows 8 & 9) with refreshed data
# - Refresh Corea del Sur (row 11) "Nuevos casos" / "Muertes hoy"
# - Reset "Nuevos casos" (col C) and "Muertes hoy" (col G) to 0 for every
#   other country row (5-186), since those countries have no new report
#   as of this data pull
# - Bump the "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 04:16"

# Blanket reset: col C (Nuevos casos) and col G (Muertes hoy) -> 0 for all
# country rows. The three special rows below are overwritten afterwards.
$ws.Range("C5:C186").Value = 0
$ws.Range("G5:G186").Value = 0

# Row 8: now Iran (refreshed counts, new-cases/deaths-today both back to 0)
$ws.Cells.Item(8, 1).Value = "Iran"
$ws.Cells.Item(8, 2).Value = 19644
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 6745
$ws.Cells.Item(8, 5).Value = 11466
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 1433

# Row 9: now Estados Unidos (refreshed counts with new cases/deaths today)
$ws.Cells.Item(9, 1).Value = "Estados Unidos"
$ws.Cells.Item(9, 2).Value = 19640
$ws.Cells.Item(9, 3).Value = 247
$ws.Cells.Item(9, 4).Value = 147
$ws.Cells.Item(9, 5).Value = 19229
$ws.Cells.Item(9, 6).Value = 64
$ws.Cells.Item(9, 7).Value = 8
$ws.Cells.Item(9, 8).Value = 264

# Row 11: Corea del Sur refreshed new-cases/deaths-today
$ws.Cells.Item(11, 3).Value = 147
$ws.Cells.Item(11, 7).Value = 8
